# Actualización automática hashcode lun jun 24 02:18:26 CEST 2019
# Update hash values (column B) for the rows whose code (column A) matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'B100' = '85819c9b0ee864700a6fb3abf7b62758'
    'B104' = 'afc45b0ea45fcd2114d8102997488408'
    'B113' = '956b266fd844e9f3fca2194ee278fadb'
    'B122' = 'd15ca3c8fb72fbbd22db7c2394f28a69'
    'B164' = '0a80cf60deec27272e68c8141fbee685'
    'B230' = 'a7ccd9496d18261177551264266f67e7'
    'B233' = '380c5e4c6ed05e85df43317f9a0cfa66'
    'B331' = 'd9986ed4380897b50d61c0803314de7c'
    'B342' = '052d5b4453144717d9154004c40aed09'
    'B343' = '9c8e173b79f48d63f00af95644862e76'
    'B605' = 'ab1ceeb631e49e4746f580da0f2b18c5'
    'B619' = 'bd09cfb4e9f5a5a1edc58ee2f6cbef23'
    'B623' = '5df9e1ffb7ca51b90d6720532ccfee6f'
    'B628' = 'ae8a27b09551a4de674da30e82a0e23c'
    'B757' = '7f5feba9a45735fef77978b4a7635326'
    'B760' = '9209ba79ce1d304c7323ecb2b6096f0b'
    'B763' = 'ee7ce75e4ca5bd240128003eac3ae32c'
    'B767' = '9cbc1dfaf8e938328934e72bd6a22359'
    'B779' = 'babf3fd530aff2ea45435a4292853ff1'
    'B818' = '4c2ed9e49577e877cba8646fab52dc00'
    'B831' = '3ebef27ff7385eb5bb0c6c1d9dc07834'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
